$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 12 ---
$ws.Range("E12").Value = 1361974149.3
$ws.Range("F12").Value = 1305498329

# --- Row 13 ---
$ws.Range("E13").Value = 332114255.61000007
$ws.Range("F13").Value = 325268233.60000002

# --- Row 14 ---
$ws.Range("E14").Value = 537985.22000000253

# --- Row 16 ---
$ws.Range("E16").Value = -60473972.810000002
$ws.Range("F16").Value = -56316441.740000002

# --- Row 18: E18 becomes a formula (was a static value) ---
$ws.Range("E18").Formula = "=SUM(E12:E17)"

# --- Row 19 ---
$ws.Range("E19").Value = -384700000.00000012

# --- Row 21: E21 becomes a formula (was a static value) ---
$ws.Range("E21").Formula = "=SUM(E18:E20)"

# --- Row 22 ---
$ws.Range("E22").Value = -20015625

# --- Row 26 ---
$ws.Range("E26").Value = 1018613403.9980445
$ws.Range("F26").Value = 1026703455

$excel.CalculateFullRebuild()
